$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 113
$ws.Range("F3").Value = 1276
$ws.Range("F4").Value = 913
$ws.Range("F5").Value = 949
$ws.Range("F6").Value = 1700
$ws.Range("F7").Value = 370
$ws.Range("F8").Value = 1139
$ws.Range("F11").Value = 102
$ws.Range("F12").Value = 256
$ws.Range("F13").Value = 27
$ws.Range("F15").Value = 626
$ws.Range("F16").Value = 126
$ws.Range("F17").Value = 80
$ws.Range("F20").Value = 318
$ws.Range("F21").Value = 89
$ws.Range("F22").Value = 640
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 623
$ws.Range("F25").Value = 124
$ws.Range("F26").Value = 28
$ws.Range("F27").Value = 833
$ws.Range("F28").Value = 291
$ws.Range("F29").Value = 86
$ws.Range("F30").Value = 21
$ws.Range("F31").Value = 242
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 6
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 23
$ws.Range("F7").Value = 233
$ws.Range("F11").Value = 106
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 294
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 294
$ws.Range("F3").Value = 113
$ws.Range("F4").Value = 1276
$ws.Range("F5").Value = 913
$ws.Range("F6").Value = 949
$ws.Range("F7").Value = 1700
$ws.Range("F8").Value = 370
$ws.Range("F9").Value = 1139
$ws.Range("F13").Value = 102
$ws.Range("F14").Value = 256
$ws.Range("F15").Value = 27
$ws.Range("F17").Value = 626
$ws.Range("F18").Value = 126
$ws.Range("F19").Value = 80
$ws.Range("F25").Value = 318
$ws.Range("F26").Value = 23
$ws.Range("F27").Value = 233
$ws.Range("F28").Value = 233
$ws.Range("F29").Value = 89
$ws.Range("F30").Value = 640
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 623
$ws.Range("F33").Value = 124
$ws.Range("F34").Value = 28
$ws.Range("F35").Value = 833
$ws.Range("F36").Value = 291
$ws.Range("F39").Value = 86
$ws.Range("F40").Value = 21
$ws.Range("F41").Value = 242
$ws.Range("F43").Value = 106
$ws.Range("F44").Value = 106
$ws.Range("F45").Value = 7
$ws.Range("F46").Value = 6
